$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 4

$ws.Cells.Item($row, 1).Value = 42607.891631944447
$ws.Cells.Item($row, 2).Value = -36
$ws.Cells.Item($row, 3).Value = 38
$ws.Cells.Item($row, 4).Value = 60
$ws.Cells.Item($row, 5).Value = 5
$ws.Cells.Item($row, 6).Value = 94
$ws.Cells.Item($row, 7).Value = 12319
$ws.Cells.Item($row, 8).Value = 21167
$ws.Cells.Item($row, 9).Value = 2275
$ws.Cells.Item($row, 10).Value = 226
$ws.Cells.Item($row, 11).Value = 357
$ws.Cells.Item($row, 12).Value = 1
$ws.Cells.Item($row, 13).Value = 18
$ws.Cells.Item($row, 14).Value = "Bag"
